$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("day sliding")
$ws.Activate()

$cols = @("D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD")

$vals117 = @(31.459894543635802,18.5937613105846,1.49942198772686,8.3354882886577499,10.168161786263299,3.63141257072634,3.3279260701063702,8.5257875628452098,2.18999681499854,0.66079255642718804,1.7120986144772301,1.04286390106311,1.24604959225402,0.27036638129784102,1.4629522778995001,9.9485649711311996,0.84205365544404998,10.096078972665699,0.84109573674310401,9.8217621380487508,7.0878702815849701,1.4627173020584801,2.2584040824045801,4.1569504118447798,6.2380145287178799,1.48117420540989,3.1914172969300001)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "117").Value = $vals117[$i]
}

$vals118 = @(31.459894543635802,11.7453225386335,2.3737019100098098,8.3354882886577499,10.168161786263299,3.63141257072634,3.3279260701063702,8.5257875628452098,2.18999681499854,0.66079255642718804,1.7120986144772301,1.04286390106311,1.24604959225402,0.27036638129784102,1.4629522778995001,9.9485649711311996,0.84205365544404998,10.096078972665699,0.84109573674310401,9.8217621380487508,7.0878702815849701,1.4627173020584801,2.2584040824045801,4.1569504118447798,6.2380145287178799,1.48117420540989,3.1914172969300001)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "118").Value = $vals118[$i]
}

$vals119 = @(29.599999999999898,17.441330224498302,1.49185868652678,2.2131735371629699,8.7255794817166201,2.6821622452719298,2.6711571208845202,9.7195304822190796,2.20097046268378,1.0507783338008601,7.4071514913996097,2.1942825153477798,0.68888319023443101,0.46348610234611498,3.0891642703740199,1.4811759324995899,2.56284891575176,8.6573606366628297,1.40042227030698,0.89499490725331599,23.930102522885502,6.8745099484647199,1.81697830331054,1.00761117840908,1.7095395146918899,3.23717818406061,3.3209584522613)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "119").Value = $vals119[$i]
}

$vals120 = @(27.799999999999901,10.285912696499,2.3546768006541199,1.65028350939709,2.09968751558011,0.90382574626335999,1.06373443281111,15.0806768864361,1.42718498561315,6.7565245602503801,7.8210289279537601,0.235385799737415,0.13344908654881901,0.66073233234818596,0.045005093083796903,0.28438825605280299,0.63050236986788799,18.0819151193047,1.13043900007137,0.371934486836066,32.721305741341197,0.0075087635519859603,0.0078999426021266403,2.9055380053674198,5.6702021231256303,0.306602093934223,0.0042452219211687597)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "120").Value = $vals120[$i]
}

$ws.Range("D117:AD120").Style = "Normal"
